$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1. Paragraph "If I were to rate this project ... rookie team at this
#    particular challenge." -- merge the grammar-checker-split runs back
#    into a clean pair of runs (tab+text, then pagebreak+text) and drop the
#    stray <w:proofErr/> markers, without disturbing the other paragraphs.
# -------------------------------------------------------------------------
$ratingPara = $d.Paragraphs.Item(5)
$ratingRange = $d.Range($ratingPara.Range.Start, $ratingPara.Range.End - 1)

$part1 = "If I were to rate this project on a scale from one to ten I would give it an eight out of ten. I chose this ratting for several reasons. Firstly, I chose to rate this project higher than average because I think it’s really freaking cool. However, this project also seemed that it may be trying to pack to many feature into too little time. Generally, it seemed this project suffered a bit from integration issues when it came to "
$part2 = "putting everything together. However it should be said that though this team has some issues they did better than any other rookie team at this particular challenge."

$ratingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
             '<w:r><w:tab/><w:t xml:space="preserve">' + $part1 + '</w:t></w:r>' +
             '<w:r><w:lastRenderedPageBreak/><w:t>' + $part2 + '</w:t></w:r>' +
             '</w:p>'
$ratingRange.InsertXML($ratingXml)

# -------------------------------------------------------------------------
# 2. Add a new paragraph right after it: "Interview With: Kevin Turkington."
# -------------------------------------------------------------------------
$afterRating = $d.Paragraphs.Item(5).Range.Duplicate
$afterRating.Collapse(0)
$afterRating.InsertParagraphAfter()
$interviewPara = $d.Paragraphs.Item(6)
$interviewPara.Range.InsertAfter("Interview With: Kevin Turkington.")

# -------------------------------------------------------------------------
# 3. Add a new paragraph right after that one: "Article By: Thomas Noelcke ."
#    -- the (formerly hidden) _GoBack bookmark ends up here, at the very
#    end of the new final content paragraph.
# -------------------------------------------------------------------------
$afterInterview = $d.Paragraphs.Item(6).Range.Duplicate
$afterInterview.Collapse(0)
$afterInterview.InsertParagraphAfter()
$bylinePara = $d.Paragraphs.Item(7)
$bylineRange = $d.Range($bylinePara.Range.Start, $bylinePara.Range.End)
$bylineXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
             '<w:r><w:t>Article By: Thomas Noelcke</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '<w:r><w:t>.</w:t></w:r>' +
             '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
             '<w:bookmarkEnd w:id="0"/>' +
             '</w:p>'
$bylineRange.InsertXML($bylineXml)

# -------------------------------------------------------------------------
# 4. Remove the old trailing paragraphs that are no longer wanted: the two
#    tab-only paragraphs, the blank paragraph, and the "Stake Holders: USLI
#    Team " paragraph -- leaving just the original two trailing empty
#    paragraphs before the section break.
# -------------------------------------------------------------------------
$deleteStart = $d.Paragraphs.Item(8).Range.Start
$deleteEnd = $d.Paragraphs.Item(11).Range.End
$deleteRange = $d.Range($deleteStart, $deleteEnd)
$deleteRange.Delete()
